# release V2021-1-2, happy new year
#
# Swap the colors of the two legend swatch rectangles (rc58 <-> rc60) and
# swap the text of the two legend-label runs (tx61 <-> tx63).
#
# All of these shapes live inside the single top-level group shape on the
# slide, so we reach them through GroupItems.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# rc58 ("F8766D" -> "619CFF"); Fill.ForeColor.RGB is read/written as 0xBBGGRR.
# Re-assert full opacity afterwards so the explicit <a:alpha val="100000"/>
# child is preserved on the solid fill, matching the original markup.
$rc58 = $grp.GroupItems.Item(56)
$rc58.Fill.ForeColor.RGB = 0xFF9C61
$rc58.Fill.Transparency = 0

# rc60 ("619CFF" -> "F8766D")
$rc60 = $grp.GroupItems.Item(58)
$rc60.Fill.ForeColor.RGB = 0x6D76F8
$rc60.Fill.Transparency = 0

# tx61 ("三" -> "一")
$tx61 = $grp.GroupItems.Item(59)
$tx61.TextFrame.TextRange.Text = "一"

# tx63 ("一" -> "三")
$tx63 = $grp.GroupItems.Item(61)
$tx63.TextFrame.TextRange.Text = "三"
